$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Folha1")

$ws.Range("A7").Value = 7
$ws.Range("A8").Value = 8

$ws.Range("A9").Select()
